$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the "Prepare"/"Assignment" headers to the new wording ---
$ws.Range("D1").Value = "Preparation need from instructor"
$ws.Range("E1").Value = "Assignment to read"

# --- Fill in the "Assignment to read" column for the remaining weeks (rows 7-14) ---
$ws.Range("E7").Value  = "chp 9, 10, 11"
$ws.Range("E8").Value  = "chp 12, 13"
$ws.Range("E9").Value  = "chp 14, 15"
$ws.Range("E10").Value = "chp 16"
$ws.Range("E11").Value = "chp 17, 18, 19"
$ws.Range("E12").Value = "chp 20 and 21"
$ws.Range("E13").Value = "chp 27"
$ws.Range("E14").Value = "Happy Git for R (Optional)"

# --- Widen columns C, D & E so the new/longer text fits nicely ---
$ws.Columns("C").ColumnWidth = 30.25
$ws.Columns("D").ColumnWidth = 43.1
$ws.Columns("E").ColumnWidth = 23.25

# --- Leave the selection on the cell below the last row, as in the saved file ---
$null = $ws.Range("E15").Select()
